# Update cryptos list (prices / 1h volume%) to the latest scraped values.
# Column D ("Price") and E ("Volume(1h)") are stored as plain text in this
# workbook (not numbers), so for any new Price value that Excel could
# otherwise auto-convert into a number we force the cell to Text format
# ("@") before assigning it. That keeps the stored value exactly as the
# original text string (e.g. "0.9992"), matching how the sheet was built.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.240.68"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "1.830.48"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.26"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07351"
$ws.Range("E8").Value = "  -1.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2932"
$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.27"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07649"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "1.832.92"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.989"
$ws.Range("E13").Value = "  -0.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6744"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.63"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008929"
$ws.Range("E16").Value = "  -2.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.876"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").Value = "29.236.86"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("D19").Value = "2.095.75"
$ws.Range("E19").Value = "  +1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.88"
$ws.Range("E20").Value = "  -1.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.52"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.392"
$ws.Range("E23").Value = "  +2.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.73"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.556"
$ws.Range("E26").Value = "  +0.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1393"
$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.64"
$ws.Range("E28").Value = "  -1.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.495"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05805"
$ws.Range("E30").Value = "  +3.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.236"
$ws.Range("E31").Value = "  +3.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.091"
$ws.Range("E32").Value = "  -0.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.105"
$ws.Range("E33").Value = "  -0.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.857"
$ws.Range("E34").Value = "  +0.87%  "

$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7217"
$ws.Range("E36").Value = "  -2.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.614"
$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.864"
$ws.Range("E38").Value = "  +3.40%  "

$ws.Range("D39").Value = "1.222.98"
$ws.Range("E39").Value = "  +0.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01763"
$ws.Range("E40").Value = "  -1.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.210"
$ws.Range("E41").Value = "  -3.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9082"
$ws.Range("E42").Value = "  +1.40%  "

$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").Value = "2.006.55"
$ws.Range("E44").Value = "  +1.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.80"
$ws.Range("E45").Value = "  +0.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.75"
$ws.Range("E46").Value = "  +0.55%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000119"
$ws.Range("E47").Value = "  -1.21%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5048"
$ws.Range("E48").Value = "  -0.71%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.185"
$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4041"
$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1172"
$ws.Range("E51").Value = "  +5.81%  "
